$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Lesson 24 (row 27): title tweak "Section Project" -> "Step Project", and
# the lesson now has a YouTube recording, so fill in the link column (F)
# with the multi-part video links - same pattern used by other multi-part
# lessons (e.g. row 19).
$ws.Range("C27").Value = "Java 8 continued, Step Project #1"
$ws.Range("F27").Value = "Part #1: https://youtu.be/seCQZfHx_bE`nPart #2: https://youtu.be/rT3auL6oukk`nPart #3: https://youtu.be/CPdlyJpnOCo"

# Match the row's look & feel to the other "recorded lesson" rows (e.g. row 19):
# centered vertical alignment, wrapped text for the name/link columns, filled
# "done" highlight on the lesson number. Column E (date) keeps its existing format.
$ws.Range("A19").Copy()
$ws.Range("A27").PasteSpecial(-4122)

$ws.Range("B19").Copy()
$ws.Range("B27").PasteSpecial(-4122)

$ws.Range("C19").Copy()
$ws.Range("C27").PasteSpecial(-4122)

$ws.Range("D19").Copy()
$ws.Range("D27").PasteSpecial(-4122)

$ws.Range("F19").Copy()
$ws.Range("F27").PasteSpecial(-4122)

# Row grows taller to fit the wrapped, multi-line link text.
$ws.Rows.Item(27).RowHeight = 41.75

# Author zoomed in a bit and moved on to working around row 33.
$excel.ActiveWindow.Zoom = 120
$ws.Range("C33").Select() | Out-Null
